$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 76

# Column A holds a date-like label ("01-04-2021") that must be stored as TEXT
# (matching the rest of the "Serie" column), not auto-converted to a date
# serial number. Writing it through a formula that yields a string, then
# copy/paste-special-values into the target cell, preserves the text type
# without Excel's literal-entry date autodetection and without minting any
# new cell style/number-format in styles.xml.
$scratch = $ws.Cells.Item(1000, 1)
$scratch.Formula = "=""01-04-2021"""
$scratch.Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4163)
$scratch.Clear()

$ws.Cells.Item($newRow, 2).Value = 26225
$ws.Cells.Item($newRow, 3).Value = 15278
$ws.Cells.Item($newRow, 4).Value = 105
$ws.Cells.Item($newRow, 5).Value = 15173
$ws.Cells.Item($newRow, 6).Value = 10865
$ws.Cells.Item($newRow, 7).Value = 7866
$ws.Cells.Item($newRow, 8).Value = 2999
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
$ws.Cells.Item($newRow, 11).Value = 0
$ws.Cells.Item($newRow, 12).Value = 82
$ws.Cells.Item($newRow, 13).Value = 0
$ws.Cells.Item($newRow, 14).Value = 26225
$ws.Cells.Item($newRow, 15).Value = 26038
$ws.Cells.Item($newRow, 16).Value = 26038
$ws.Cells.Item($newRow, 17).Value = 187
